$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the three runs of "Faire " + "l'achat de la batterie" + " " into a
#    single run "Faire l'achat de la batterie " (paragraph with that text).
# ---------------------------------------------------------------------------
$rngBattery = $d.Content
$found = $rngBattery.Find.Execute("Faire l" + [char]0x2019 + "achât de la batterie", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $p = $rngBattery.Paragraphs(1)
    $pr = $p.Range
    $pr.MoveEnd(1, -1)
    # Force a real text mutation so the engine re-emits a single run instead
    # of leaving the original run split untouched.
    $pr.Text = ""
    $pr.Text = "Faire l" + [char]0x2019 + "achât de la batterie "
}

# ---------------------------------------------------------------------------
# 2. Split "Mettre en évidence les critères pour chaque livrable" so that the
#    word "critères" becomes "activités" in its own run, matching the
#    three-run structure produced by a normal Word "type-over-selection"
#    edit: "Mettre en évidence les " / "activités" / " pour chaque livrable".
# ---------------------------------------------------------------------------
$rngCriteres = $d.Content
$found = $rngCriteres.Find.Execute("critères pour chaque livrable", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $wordRng = $d.Content
    $wordRng.Find.Execute("critères", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $wordRng.Text = "activités"
    # Toggling Bold on/off forces the new run to stay distinct from its
    # neighbours (which share identical formatting) without leaving any
    # residual explicit formatting behind.
    $wordRng.Bold = 1
    $wordRng.Bold = 0
}

# ---------------------------------------------------------------------------
# 3. Remove the old _GoBack bookmark that used to sit right after
#    "Installer la batterie".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4. Re-create a fresh _GoBack bookmark (plus a new _Hlk61104030 bookmark)
#    spanning the title paragraphs "ROVUS : Structure de découpage" and
#    "Module de puissance V1".
# ---------------------------------------------------------------------------
$pTitle1 = $d.Paragraphs(19)   # "ROVUS : Structure de découpage"
$pTitle2 = $d.Paragraphs(20)   # "Module de puissance V1"
$titleStart = $pTitle1.Range.Start
$titleEnd = $pTitle2.Range.End
$titleRange = $d.Range($titleStart, $titleEnd)
$d.Bookmarks.Add("_Hlk61104030", $titleRange)
$titleRange2 = $d.Range($titleStart, $titleEnd)
$d.Bookmarks.Add("_GoBack", $titleRange2)

# ---------------------------------------------------------------------------
# 5. Add a new _Hlk61103909 bookmark spanning the "Édouard Villemure" and
#    "2020-12-30" paragraphs.
# ---------------------------------------------------------------------------
$pAuthor = $d.Paragraphs(30)   # "Édouard Villemure"
$pDate = $d.Paragraphs(31)     # "2020-12-30"
$authorStart = $pAuthor.Range.Start
$authorEnd = $pDate.Range.End
$authorRange = $d.Range($authorStart, $authorEnd)
$d.Bookmarks.Add("_Hlk61103909", $authorRange)

# ---------------------------------------------------------------------------
# 6. Add a new _Hlk61103673 bookmark spanning from "Choisir une batterie"
#    through to the empty Heading 3 paragraph at the end of the document.
# ---------------------------------------------------------------------------
$pChoisir = $d.Paragraphs(32)     # "Choisir une batterie"
$pHeading3 = $d.Paragraphs(134)   # trailing empty Heading 3 paragraph
$choisirStart = $pChoisir.Range.Start
$choisirEnd = $pHeading3.Range.End
$choisirRange = $d.Range($choisirStart, $choisirEnd)
$d.Bookmarks.Add("_Hlk61103673", $choisirRange)

Write-Output "done"
